$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F ("dSF")
$changes = @{
    2  = -5
    3  = -5
    5  = 6
    6  = -5
    11 = 6
    12 = 1
    13 = -4
    15 = 2
    21 = -2
    23 = 3
    27 = -1
    28 = 1
    31 = 6
    37 = -4
    38 = 0
    41 = 2
    44 = -2
    46 = 1
    49 = -3
    50 = -9
    51 = 5
    54 = -9
    56 = -6
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
